# Updates cryptos list: Price (D) and Volume(1h) (E) columns per diff
# (mirrors the upstream GitHub Actions scraper refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.788.45"
$ws.Range("E2").Value = "  +2.79%  "
$ws.Range("D3").Value = "1.722.98"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("E4").Value = "  -0.73%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.54"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("E8").Value = "  +12.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.266"
$ws.Range("E9").Value = "  +4.64%  "
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "1.965.92"
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("D13").Value = "1.730.59"
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.566"
$ws.Range("E15").Value = "  +6.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.01"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").Value = "27.793.09"
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.62"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.01"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("D20").Value = "0.0₃0749"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.63"
$ws.Range("E22").Value = "  +3.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.74"
$ws.Range("E23").Value = "  +5.33%  "
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.57"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.57"
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.62"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").Value = "1.551.58"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("E34").Value = "  +4.43%  "
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.970"
$ws.Range("E36").Value = "  +6.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.617"
$ws.Range("E37").Value = "  +4.67%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.55"
$ws.Range("E41").Value = "  +5.37%  "
$ws.Range("E42").Value = "  +5.79%  "
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").Value = "1.867.90"
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.793"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("E47").Value = "  +9.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "92.22"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").Value = "0.0₆0111"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("E50").Value = "  +3.98%  "
$ws.Range("E51").Value = "  +2.36%  "
